$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("mean")
$ws2 = $wb.Worksheets.Item("stdev")
$ws3 = $wb.Worksheets.Item("summary")

# sheet "mean": update H2O (K) and N (L) columns, rows 2-10
$ws1.Range("K2").Value = 11.864
$ws1.Range("L2").Value = 0.556
$ws1.Range("K3").Value = 11.471
$ws1.Range("L3").Value = 0.948
$ws1.Range("K4").Value = 11.277
$ws1.Range("L4").Value = 1.143
$ws1.Range("K5").Value = 13.211
$ws1.Range("L5").Value = 0.219
$ws1.Range("K6").Value = 12.828
$ws1.Range("L6").Value = 0.602
$ws1.Range("K7").Value = 12.705
$ws1.Range("L7").Value = 0.725
$ws1.Range("K8").Value = 7.82
$ws1.Range("L8").Value = 0.649
$ws1.Range("K9").Value = 7.417
$ws1.Range("L9").Value = 1.052
$ws1.Range("K10").Value = 7.202
$ws1.Range("L10").Value = 1.267

# sheet "stdev": update H2O_sd (K) and N_sd (L) columns, rows 2-10
$ws2.Range("K2").Value = 0.034
$ws2.Range("L2").Value = 0.034
$ws2.Range("K3").Value = 0.037
$ws2.Range("L3").Value = 0.037
$ws2.Range("K4").Value = 0.044
$ws2.Range("L4").Value = 0.044
$ws2.Range("K5").Value = 0.108
$ws2.Range("L5").Value = 0.108
$ws2.Range("K6").Value = 0.105
$ws2.Range("L6").Value = 0.105
$ws2.Range("K7").Value = 0.126
$ws2.Range("L7").Value = 0.126
$ws2.Range("K8").Value = 0.03
$ws2.Range("L8").Value = 0.03
$ws2.Range("K9").Value = 0.024
$ws2.Range("L9").Value = 0.024
$ws2.Range("K10").Value = 0.029
$ws2.Range("L10").Value = 0.029

# sheet "summary": update H2O (H), H2O_sd (I), N (P), N_sd (Q) columns, rows 2-10
$ws3.Range("H2").Value = 11.864
$ws3.Range("I2").Value = 0.034
$ws3.Range("P2").Value = 0.556
$ws3.Range("Q2").Value = 0.034
$ws3.Range("H3").Value = 11.471
$ws3.Range("I3").Value = 0.037
$ws3.Range("P3").Value = 0.948
$ws3.Range("Q3").Value = 0.037
$ws3.Range("H4").Value = 11.277
$ws3.Range("I4").Value = 0.044
$ws3.Range("P4").Value = 1.143
$ws3.Range("Q4").Value = 0.044
$ws3.Range("H5").Value = 13.211
$ws3.Range("I5").Value = 0.108
$ws3.Range("P5").Value = 0.219
$ws3.Range("Q5").Value = 0.108
$ws3.Range("H6").Value = 12.828
$ws3.Range("I6").Value = 0.105
$ws3.Range("P6").Value = 0.602
$ws3.Range("Q6").Value = 0.105
$ws3.Range("H7").Value = 12.705
$ws3.Range("I7").Value = 0.126
$ws3.Range("P7").Value = 0.725
$ws3.Range("Q7").Value = 0.126
$ws3.Range("H8").Value = 7.82
$ws3.Range("I8").Value = 0.03
$ws3.Range("P8").Value = 0.649
$ws3.Range("Q8").Value = 0.03
$ws3.Range("H9").Value = 7.417
$ws3.Range("I9").Value = 0.024
$ws3.Range("P9").Value = 1.052
$ws3.Range("Q9").Value = 0.024
$ws3.Range("H10").Value = 7.202
$ws3.Range("I10").Value = 0.029
$ws3.Range("P10").Value = 1.267
$ws3.Range("Q10").Value = 0.029

Write-Output "Updated H2O and N columns on mean, stdev, and summary sheets."
